$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F, reusing the same formatting as the other
# header cells (bold / bordered / centered) by copying E1's format.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$timestamps = @(
    "2021-10-05 10:52:15.612165",
    "2021-10-05 10:52:15.612178",
    "2021-10-05 10:52:15.612182",
    "2021-10-05 10:52:15.612185",
    "2021-10-05 10:52:15.612189",
    "2021-10-05 10:52:15.612192",
    "2021-10-05 10:52:15.612195",
    "2021-10-05 10:52:15.612198",
    "2021-10-05 10:52:15.612201",
    "2021-10-05 10:52:15.612204",
    "2021-10-05 10:52:15.612208",
    "2021-10-05 10:52:15.612211",
    "2021-10-05 10:52:15.612214",
    "2021-10-05 10:52:15.612217",
    "2021-10-05 10:52:15.612219",
    "2021-10-05 10:52:15.612228",
    "2021-10-05 10:52:15.612231",
    "2021-10-05 10:52:15.612234",
    "2021-10-05 10:52:15.612237",
    "2021-10-05 10:52:15.612256",
    "2021-10-05 10:52:15.612258",
    "2021-10-05 10:52:15.612261",
    "2021-10-05 10:52:15.612263",
    "2021-10-05 10:52:15.612266",
    "2021-10-05 10:52:15.612269",
    "2021-10-05 10:52:15.612271",
    "2021-10-05 10:52:15.612274",
    "2021-10-05 10:52:15.612277",
    "2021-10-05 10:52:15.612279",
    "2021-10-05 10:52:15.612282",
    "2021-10-05 10:52:15.612284",
    "2021-10-05 10:52:15.612287",
    "2021-10-05 10:52:15.612292",
    "2021-10-05 10:52:15.612295",
    "2021-10-05 10:52:15.612297",
    "2021-10-05 10:52:15.612300",
    "2021-10-05 10:52:15.612303",
    "2021-10-05 10:52:15.612305",
    "2021-10-05 10:52:15.612308",
    "2021-10-05 10:52:15.612310",
    "2021-10-05 10:52:15.612314",
    "2021-10-05 10:52:15.612316",
    "2021-10-05 10:52:15.612319",
    "2021-10-05 10:52:15.612321",
    "2021-10-05 10:52:15.612324",
    "2021-10-05 10:52:15.612326"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
